$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = $true
$ws.Range("D2").Value = $true
$ws.Range("D2").Select()
